$d = $word.ActiveDocument

function Replace-RangeWithXml($range, $innerXml) {
    # NOTE: a Range produced by Find.Execute() (even via .Duplicate()) only
    # *inserts* XML at its Start when handed to InsertXML; re-materialising
    # the same Start/End via $d.Range(...) gives a plain Range that properly
    # replaces its span instead.
    $fresh = $d.Range($range.Start, $range.End)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $innerXml +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $fresh.InsertXML($pkg)
}

# ------------------------------------------------------------------
# 1) "Add some comments ... GitHub user id ..." paragraph:
#    split the single run so "GitHub" is wrapped in a spellStart/spellEnd
#    proofErr pair (as Word's background spell-checker would do).
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Add some comments about Version management after this line, or just add some text so there is a change to this file.  Remember that your GitHub user id must be submitted in you assignment report!", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$xml1 = '<w:p>' +
        '<w:r><w:t xml:space="preserve">Add some comments about Version management after this line, or just add some text so there is a change to this file.  Remember that your </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>GitHub</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> user id must be submitted in you assignment report!</w:t></w:r>' +
        '</w:p>'
Replace-RangeWithXml $r1 $xml1

# ------------------------------------------------------------------
# 2) ">>>  your stuff after this line >>>" paragraph merges with the
#    following "Ths is a test of changes" paragraph: the first gets
#    split into 3 runs with a gramStart/gramEnd proofErr pair, and the
#    _GoBack bookmark (previously on the second paragraph) is kept,
#    now trailing the rebuilt run; the second paragraph's own text is
#    dropped entirely (its paragraph mark disappears too).
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute(">>>  your stuff after this line >>>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$paraIndex = $r2.Paragraphs(1).Index
$p1 = $d.Paragraphs($paraIndex)
$p2 = $d.Paragraphs($paraIndex + 1)
$mergedRange = $d.Range($p1.Range.Start, $p2.Range.End)

$xml2 = '<w:p>' +
        '<w:r><w:t>&gt;&gt;</w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>&gt;  your</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t xml:space="preserve"> stuff after this line &gt;&gt;&gt;</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '</w:p>'
Replace-RangeWithXml $mergedRange $xml2

# ------------------------------------------------------------------
# 3) "Baz changes" paragraph becomes "This is a test of changes made by
#    Crispin", split across four runs ("Th" / "i" / "s is a test of
#    changes" / " made by Crispin").
# ------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("Baz changes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$xml3 = '<w:p>' +
        '<w:r><w:t>Th</w:t></w:r>' +
        '<w:r><w:t>i</w:t></w:r>' +
        '<w:r><w:t>s is a test of changes</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> made by Crispin</w:t></w:r>' +
        '</w:p>'
Replace-RangeWithXml $r3 $xml3
